$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 45250

$ws.Range("E15").Value = "Poursuite du CDC, charte graph, moodboard, arbo"

$ws.Range("E26").Select()
